$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.358.50"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.934.78"
$ws.Range("E3").Value = "  +1.33%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.80"
$ws.Range("E5").Value = "  +2.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7205"
$ws.Range("E6").Value = "  +4.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3283"
$ws.Range("E8").Value = "  +2.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.54"
$ws.Range("E9").Value = "  +8.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07238"
$ws.Range("E10").Value = "  +6.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8044"
$ws.Range("E11").Value = "  +2.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08090"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.932.79"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.432"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.62"
$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.02"
$ws.Range("E16").Value = "  +5.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.342.09"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.44"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008211"
$ws.Range("E19").Value = "  +5.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.805"
$ws.Range("E20").Value = "  +0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.187.15"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.935"
$ws.Range("E24").Value = "  +2.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.712"
$ws.Range("E25").Value = "  +2.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.10"
$ws.Range("E26").Value = "  +4.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.343"
$ws.Range("E27").Value = "  +7.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.28"
$ws.Range("E28").Value = "  +3.70%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.546"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.437"
$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("E33").Value = "  +1.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05228"
$ws.Range("E34").Value = "  +4.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.264"
$ws.Range("E35").Value = "  +7.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7488"
$ws.Range("E36").Value = "  +1.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.764"
$ws.Range("E37").Value = "  +1.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01966"
$ws.Range("E38").Value = "  +3.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.803"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.09"
$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.446"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4535"
$ws.Range("E42").Value = "  +3.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.030"
$ws.Range("E43").Value = "  +2.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8427"
$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.88"
$ws.Range("E46").Value = "  +0.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.810"
$ws.Range("E47").Value = "  +1.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.442"
$ws.Range("E48").Value = "  +4.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.78"
$ws.Range("E49").Value = "  +3.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4185"
$ws.Range("E50").Value = "  +4.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06057"
$ws.Range("E51").Value = "  +2.68%  "
